# Update the "想去人数" (wanted-to-go count) figures in F3, F7, F12
# on both the "展览" sheet and the "全部类型" sheet (which mirrors it).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 6398
    $ws.Range("F7").Value = 1918
    $ws.Range("F12").Value = 5595
}
